$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to coerce numeric-looking tracking numbers into text
# (matches the workbook's existing convention of storing these IDs as
# shared strings rather than numbers) without leaving any NumberFormat /
# style residue on the target cell: build the literal text via TEXT(),
# copy it, then paste *values only* into the destination.
$scratch = $ws.Range("ZZ1000")

function Set-TextValue($cellref, $val) {
    $scratch.Formula = "=TEXT(" + $val + ",""0"")"
    $scratch.Copy()
    $ws.Range($cellref).PasteSpecial(-4163)
}

Set-TextValue "C2" "320017958830"
Set-TextValue "C3" "320017958841"
Set-TextValue "C4" "320017958874"
Set-TextValue "C5" "320017958896"
Set-TextValue "D5" "320017958896"
Set-TextValue "C6" "320017958933"
Set-TextValue "D6" "320017958933"
Set-TextValue "C7" "320017958955"
Set-TextValue "D7" "320017958955"
Set-TextValue "C8" "320017958988"
Set-TextValue "C9" "320017959002"
Set-TextValue "C10" "320017959035"
Set-TextValue "C11" "320017959057"
Set-TextValue "C12" "320017959090"
Set-TextValue "C13" "320017959116"
Set-TextValue "D13" "320017959116"
Set-TextValue "C14" "320017959149"
Set-TextValue "D14" "320017959149"
Set-TextValue "C15" "320017959160"
Set-TextValue "D15" "320017959160"
Set-TextValue "C16" "320017959208"
Set-TextValue "D16" "320017959208"
Set-TextValue "C17" "320017959220"
Set-TextValue "D17" "320017959220"
Set-TextValue "C18" "320017959263"
Set-TextValue "C19" "320017959285"
Set-TextValue "C20" "320017959311"
Set-TextValue "C21" "320017959333"
Set-TextValue "C22" "320017959366"

# Plain (non-numeric) text value - no coercion trick needed.
$ws.Range("Q3").Value = "Pass"

# Clean up the scratch cell so it doesn't leave any trace in the sheet.
$scratch.Clear()
